# "Natmi following Dr Hou advice"
# Recomputed Cxcl10/Cxcr3 ligand-receptor stats: sending/target cluster counts
# are now aggregated per side (3 cells each) instead of per sending/target
# pair (previously 1), which changes every derived expression/specificity
# value in rows 2-5 and introduces the two additional sending clusters
# (M2, sCs) x two target clusters (ECs, M2) combinations as new rows 6-9.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (updates)
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,7).Value = 18.76479233333333
$ws.Cells.Item(2,8).Value = 56.294377
$ws.Cells.Item(2,9).Value = 0.09818846546758657
$ws.Cells.Item(2,10).Value = 0.09818846546758656
$ws.Cells.Item(2,12).Value = 0.3333333333333333
$ws.Cells.Item(2,13).Value = 1.863797
$ws.Cells.Item(2,14).Value = 5.591391
$ws.Cells.Item(2,15).Value = 0.5455000708290748
$ws.Cells.Item(2,16).Value = 0.5455000708290748
$ws.Cells.Item(2,17).Value = 34.97376365648967
$ws.Cells.Item(2,18).Value = 314.763872908407
$ws.Cells.Item(2,19).Value = 0.05356181486716664
$ws.Cells.Item(2,20).Value = 0.05356181486716663

# Row 3 (updates)
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,7).Value = 18.76479233333333
$ws.Cells.Item(3,8).Value = 56.294377
$ws.Cells.Item(3,9).Value = 0.09818846546758657
$ws.Cells.Item(3,10).Value = 0.09818846546758656
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,13).Value = 1.552879
$ws.Cells.Item(3,14).Value = 4.658637000000001
$ws.Cells.Item(3,15).Value = 0.4544999291709252
$ws.Cells.Item(3,16).Value = 0.4544999291709252
$ws.Cells.Item(3,17).Value = 29.13945195379434
$ws.Cells.Item(3,18).Value = 262.255067584149
$ws.Cells.Item(3,19).Value = 0.04462665060041993
$ws.Cells.Item(3,20).Value = 0.04462665060041993

# Row 4 (updates)
$ws.Cells.Item(4,1).Value = "FAPs"
$ws.Cells.Item(4,4).Value = "ECs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,7).Value = 72.11798333333333
$ws.Cells.Item(4,8).Value = 216.35395
$ws.Cells.Item(4,9).Value = 0.3773638413007209
$ws.Cells.Item(4,10).Value = 0.3773638413007209
$ws.Cells.Item(4,12).Value = 0.3333333333333333
$ws.Cells.Item(4,13).Value = 1.863797
$ws.Cells.Item(4,14).Value = 5.591391
$ws.Cells.Item(4,15).Value = 0.5455000708290748
$ws.Cells.Item(4,16).Value = 0.5455000708290748
$ws.Cells.Item(4,17).Value = 134.4132809827167
$ws.Cells.Item(4,18).Value = 1209.71952884445
$ws.Cells.Item(4,19).Value = 0.205852002157875
$ws.Cells.Item(4,20).Value = 0.205852002157875

# Row 5 (updates)
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,7).Value = 72.11798333333333
$ws.Cells.Item(5,8).Value = 216.35395
$ws.Cells.Item(5,9).Value = 0.3773638413007209
$ws.Cells.Item(5,10).Value = 0.3773638413007209
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,13).Value = 1.552879
$ws.Cells.Item(5,14).Value = 4.658637000000001
$ws.Cells.Item(5,15).Value = 0.4544999291709252
$ws.Cells.Item(5,16).Value = 0.4544999291709252
$ws.Cells.Item(5,17).Value = 111.9905018406833
$ws.Cells.Item(5,18).Value = 1007.91451656615
$ws.Cells.Item(5,19).Value = 0.1715118391428459
$ws.Cells.Item(5,20).Value = 0.1715118391428459

# Row 6 (new row)
$ws.Cells.Item(6,1).Value = "M2"
$ws.Cells.Item(6,2).Value = "Cxcl10"
$ws.Cells.Item(6,3).Value = "Cxcr3"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 83.31930033333333
$ws.Cells.Item(6,8).Value = 249.957901
$ws.Cells.Item(6,9).Value = 0.4359757410707099
$ws.Cells.Item(6,10).Value = 0.4359757410707098
$ws.Cells.Item(6,11).Value = 1
$ws.Cells.Item(6,12).Value = 0.3333333333333333
$ws.Cells.Item(6,13).Value = 1.863797
$ws.Cells.Item(6,14).Value = 5.591391
$ws.Cells.Item(6,15).Value = 0.5455000708290748
$ws.Cells.Item(6,16).Value = 0.5455000708290748
$ws.Cells.Item(6,17).Value = 155.2902620033657
$ws.Cells.Item(6,18).Value = 1397.612358030291
$ws.Cells.Item(6,19).Value = 0.2378247976338306
$ws.Cells.Item(6,20).Value = 0.2378247976338306

# Row 7 (new row)
$ws.Cells.Item(7,1).Value = "M2"
$ws.Cells.Item(7,2).Value = "Cxcl10"
$ws.Cells.Item(7,3).Value = "Cxcr3"
$ws.Cells.Item(7,4).Value = "M2"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 83.31930033333333
$ws.Cells.Item(7,8).Value = 249.957901
$ws.Cells.Item(7,9).Value = 0.4359757410707099
$ws.Cells.Item(7,10).Value = 0.4359757410707098
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 1.552879
$ws.Cells.Item(7,14).Value = 4.658637000000001
$ws.Cells.Item(7,15).Value = 0.4544999291709252
$ws.Cells.Item(7,16).Value = 0.4544999291709252
$ws.Cells.Item(7,17).Value = 129.3847917823263
$ws.Cells.Item(7,18).Value = 1164.463126040937
$ws.Cells.Item(7,19).Value = 0.1981509434368793
$ws.Cells.Item(7,20).Value = 0.1981509434368793

# Row 8 (new row)
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Cxcl10"
$ws.Cells.Item(8,3).Value = "Cxcr3"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 16.90787
$ws.Cells.Item(8,8).Value = 50.72361
$ws.Cells.Item(8,9).Value = 0.08847195216098278
$ws.Cells.Item(8,10).Value = 0.08847195216098278
$ws.Cells.Item(8,11).Value = 1
$ws.Cells.Item(8,12).Value = 0.3333333333333333
$ws.Cells.Item(8,13).Value = 1.863797
$ws.Cells.Item(8,14).Value = 5.591391
$ws.Cells.Item(8,15).Value = 0.5455000708290748
$ws.Cells.Item(8,16).Value = 0.5455000708290748
$ws.Cells.Item(8,17).Value = 31.51283738239
$ws.Cells.Item(8,18).Value = 283.61553644151
$ws.Cells.Item(8,19).Value = 0.04826145617020262
$ws.Cells.Item(8,20).Value = 0.04826145617020262

# Row 9 (new row)
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Cxcl10"
$ws.Cells.Item(9,3).Value = "Cxcr3"
$ws.Cells.Item(9,4).Value = "M2"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 16.90787
$ws.Cells.Item(9,8).Value = 50.72361
$ws.Cells.Item(9,9).Value = 0.08847195216098278
$ws.Cells.Item(9,10).Value = 0.08847195216098278
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 1.552879
$ws.Cells.Item(9,14).Value = 4.658637000000001
$ws.Cells.Item(9,15).Value = 0.4544999291709252
$ws.Cells.Item(9,16).Value = 0.4544999291709252
$ws.Cells.Item(9,17).Value = 26.25587625773
$ws.Cells.Item(9,18).Value = 236.30288631957
$ws.Cells.Item(9,19).Value = 0.04021049599078016
$ws.Cells.Item(9,20).Value = 0.04021049599078016
